# Update the "Local Computer" comparison table (rows 10-13) with new
# MSI output / comparison values, then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (lm): Original 2.97 -> 3.28, Parallel 3.11 -> 3.04
$ws.Range("F10").Value = 3.28
$ws.Range("G10").Value = 3.04

# Row 11 (glmnet): Original 5.78 -> 5.95, Parallel stays 1.91
$ws.Range("F11").Value = 5.95

# Row 12 (ranger): Original 36.5 -> 37.8, Parallel 35.4 -> 34.6
$ws.Range("F12").Value = 37.8
$ws.Range("G12").Value = 34.6

# Row 13 (xgbTree): Original 128 -> 129, Parallel stays 28
$ws.Range("F13").Value = 129

# Move the selection to E21 to match the author's final cursor position
$ws.Range("E21").Select()
